$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in Thursday (row 18) hours for week 3.2 section
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 4
$ws.Range("I18").Value = 4

# Update the selected cell in the sheet view
$ws.Range("J18").Select()
